$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(8,1).NumberFormat = "0.00%"
Write-Host $ws.Cells.Item(8,1).NumberFormat
